# NYPD CompStat weekly report refresh: new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume/issue number and the reporting week dates.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# ---------------------------------------------------------------------------
# Column H got narrower (matches the other "% Chg" columns' width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 5.43

# ---------------------------------------------------------------------------
# Helper pattern used throughout: when a numeric cell becomes a "no data"
# placeholder (shown as "0" / "***.*" text) or vice versa, set the raw
# value first, then copy number-format/style from a neighboring cell that
# already carries the right style so the cell's style index matches.
# ---------------------------------------------------------------------------

# Row 14 - Murder
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("M14").Value = -16.666666666666

# Row 15 - Rape
$ws.Range("M15").Value = -23.076923076923

# Row 16 - Robbery
$ws.Range("C16").Value = 3
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 158
$ws.Range("J16").Value = 171
$ws.Range("K16").Value = -7.602339181286
$ws.Range("L16").Value = -3.658536585365
$ws.Range("M16").Value = -33.891213389121
$ws.Range("N16").Value = -81.301775147929

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = -23.076923076923
$ws.Range("I17").Value = 279
$ws.Range("J17").Value = 230
$ws.Range("K17").Value = 21.304347826087
$ws.Range("L17").Value = 28.571428571428
$ws.Range("M17").Value = 61.271676300578
$ws.Range("N17").Value = -60.199714693295

# Row 18 - Burglary
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 89
$ws.Range("J18").Value = 104
$ws.Range("K18").Value = -14.423076923076
$ws.Range("L18").Value = -47.647058823529
$ws.Range("M18").Value = -13.592233009708
$ws.Range("N18").Value = -89.517078916372

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 476
$ws.Range("J19").Value = 377
$ws.Range("K19").Value = 26.259946949602
$ws.Range("L19").Value = 23.316062176165
$ws.Range("M19").Value = 61.355932203389
$ws.Range("N19").Value = 1.061571125265

# Row 20 - G.L.A.
$ws.Range("C20").Value = "'0"
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 60
$ws.Range("K20").Value = -10.447761194029
$ws.Range("L20").Value = 17.647058823529
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -57.746478873239

# Row 21 - TOTAL
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -35.714285714285
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -3.076923076923
$ws.Range("I21").Value = 1077
$ws.Range("J21").Value = 960
$ws.Range("K21").Value = 12.1875
$ws.Range("L21").Value = 7.592407592407
$ws.Range("M21").Value = 26.855123674911
$ws.Range("N21").Value = -65.066493674991

# Row 22 - Transit
$ws.Range("I22").Value = 22
$ws.Range("K22").Value = -4.347826086956
$ws.Range("L22").Value = -29.032258064516
$ws.Range("M22").Value = 57.142857142857

# Row 23 - Housing
$ws.Range("C23").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 82
$ws.Range("K23").Value = -8.536585365853

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 58.333333333333
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 58.441558441558
$ws.Range("I24").Value = 1264
$ws.Range("J24").Value = 1199
$ws.Range("K24").Value = 5.421184320266
$ws.Range("L24").Value = 11.365638766519
$ws.Range("M24").Value = 29.774127310061

# Row 25 - Retail Theft
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 159.259259259259
$ws.Range("I25").Value = 762
$ws.Range("J25").Value = 681
$ws.Range("K25").Value = 11.894273127753
$ws.Range("L25").Value = 13.392857142857

# Row 26 - Misd. Assault
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 125
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 77.777777777777
$ws.Range("I26").Value = 434
$ws.Range("J26").Value = 382
$ws.Range("K26").Value = 13.612565445026
$ws.Range("L26").Value = 10.152284263959
$ws.Range("M26").Value = -11.065573770491

# Row 28 - Other Sex Crimes
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 43
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = 16.216216216216
$ws.Range("L28").Value = -8.510638297872

# Row 29 - Shooting Vic.
$ws.Range("G29").Value = "'0"
$ws.Range("F29").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = "'***.*"
$ws.Range("E29").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("M29").Value = -17.391304347826

# Row 30 - Shooting Inc.
$ws.Range("G30").Value = "'0"
$ws.Range("F30").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = "'***.*"
$ws.Range("E30").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("M30").Value = -30
